# Auto-generated Excel COM-interop script
# Adds new scrim rows to 4 worksheets, matching the target diff.

$wb = $excel.ActiveWorkbook

# ---- Triple Dribble (sheet index 1) ----
$ws = $wb.Worksheets.Item(1)

# Copy "Equipo 2" row formatting (style pattern from existing row 4)
$ws.Range("A4:N4").Copy()
$ws.Range("A39:N39").PasteSpecial(-4122)
$ws.Range("A40:N40").PasteSpecial(-4122)

# Set cell values for new rows
$ws.Cells.Item(39,1).Value = "MEEPLE"
$ws.Cells.Item(39,2).Value = "LOU"
$ws.Cells.Item(39,3).Value = "SHADE"
$ws.Cells.Item(39,4).Value = "CROW"
$ws.Cells.Item(39,5).Value = "CORDELIUS"
$ws.Cells.Item(39,6).Value = "BARLEY"
$ws.Cells.Item(39,7).Value = "Equipo 2"
$ws.Cells.Item(39,8).Value = "IC|Nob"
$ws.Cells.Item(39,9).Value = "IC|Mebius"
$ws.Cells.Item(39,10).Value = "IC|RamaZR"
$ws.Cells.Item(39,11).Value = "Enraged 💔"
$ws.Cells.Item(39,12).Value = "SUP|Filippo神"
$ws.Cells.Item(39,13).Value = "SUP|Tomzy"
$ws.Cells.Item(39,14).Value = "20250723T162138.000Z"

$ws.Cells.Item(40,1).Value = "MEEPLE"
$ws.Cells.Item(40,2).Value = "LOU"
$ws.Cells.Item(40,3).Value = "SHADE"
$ws.Cells.Item(40,4).Value = "CROW"
$ws.Cells.Item(40,5).Value = "CORDELIUS"
$ws.Cells.Item(40,6).Value = "BARLEY"
$ws.Cells.Item(40,7).Value = "Equipo 2"
$ws.Cells.Item(40,8).Value = "IC|Nob"
$ws.Cells.Item(40,9).Value = "IC|Mebius"
$ws.Cells.Item(40,10).Value = "IC|RamaZR"
$ws.Cells.Item(40,11).Value = "Enraged 💔"
$ws.Cells.Item(40,12).Value = "SUP|Filippo神"
$ws.Cells.Item(40,13).Value = "SUP|Tomzy"
$ws.Cells.Item(40,14).Value = "20250723T162047.000Z"

# ---- Ring of Fire (sheet index 15) ----
$ws = $wb.Worksheets.Item(15)

# Copy "Equipo 2" row formatting (style pattern from existing row 4)
$ws.Range("A4:N4").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)
$ws.Range("A7:N7").PasteSpecial(-4122)

# Copy "Equipo 1" row formatting from Triple Dribble row 5 (no local template)
$wsSrc = $wb.Worksheets.Item(1)
$wsSrc.Range("A5:N5").Copy()
$ws.Range("A8:N8").PasteSpecial(-4122)
$ws.Range("A9:N9").PasteSpecial(-4122)

# Set cell values for new rows
$ws.Cells.Item(6,1).Value = "MEG"
$ws.Cells.Item(6,2).Value = "DRACO"
$ws.Cells.Item(6,3).Value = "BEA"
$ws.Cells.Item(6,4).Value = "DOUG"
$ws.Cells.Item(6,5).Value = "AMBER"
$ws.Cells.Item(6,6).Value = "GRAY"
$ws.Cells.Item(6,7).Value = "Equipo 2"
$ws.Cells.Item(6,8).Value = "NXT|Rup"
$ws.Cells.Item(6,9).Value = "NXT|Arthur"
$ws.Cells.Item(6,10).Value = "NXT|amos"
$ws.Cells.Item(6,11).Value = "FUT|Nowy297"
$ws.Cells.Item(6,12).Value = "FUT|MeOw"
$ws.Cells.Item(6,13).Value = "FUT|GeRo"
$ws.Cells.Item(6,14).Value = "20250723T162258.000Z"

$ws.Cells.Item(7,1).Value = "MEG"
$ws.Cells.Item(7,2).Value = "DRACO"
$ws.Cells.Item(7,3).Value = "BEA"
$ws.Cells.Item(7,4).Value = "DOUG"
$ws.Cells.Item(7,5).Value = "AMBER"
$ws.Cells.Item(7,6).Value = "GRAY"
$ws.Cells.Item(7,7).Value = "Equipo 2"
$ws.Cells.Item(7,8).Value = "NXT|Rup"
$ws.Cells.Item(7,9).Value = "NXT|Arthur"
$ws.Cells.Item(7,10).Value = "NXT|amos"
$ws.Cells.Item(7,11).Value = "FUT|Nowy297"
$ws.Cells.Item(7,12).Value = "FUT|MeOw"
$ws.Cells.Item(7,13).Value = "FUT|GeRo"
$ws.Cells.Item(7,14).Value = "20250723T162056.000Z"

$ws.Cells.Item(8,1).Value = "COLETTE"
$ws.Cells.Item(8,2).Value = "ALLI"
$ws.Cells.Item(8,3).Value = "MR. P"
$ws.Cells.Item(8,4).Value = "LUMI"
$ws.Cells.Item(8,5).Value = "KAZE"
$ws.Cells.Item(8,6).Value = "HANK"
$ws.Cells.Item(8,7).Value = "Equipo 1"
$ws.Cells.Item(8,8).Value = "NXT|Rup"
$ws.Cells.Item(8,9).Value = "NXT|Arthur"
$ws.Cells.Item(8,10).Value = "NXT|amos"
$ws.Cells.Item(8,11).Value = "FUT|GeRo"
$ws.Cells.Item(8,12).Value = "FUT|Nowy297"
$ws.Cells.Item(8,13).Value = "FUT|MeOw"
$ws.Cells.Item(8,14).Value = "20250723T161357.000Z"

$ws.Cells.Item(9,1).Value = "COLETTE"
$ws.Cells.Item(9,2).Value = "ALLI"
$ws.Cells.Item(9,3).Value = "MR. P"
$ws.Cells.Item(9,4).Value = "LUMI"
$ws.Cells.Item(9,5).Value = "KAZE"
$ws.Cells.Item(9,6).Value = "HANK"
$ws.Cells.Item(9,7).Value = "Equipo 1"
$ws.Cells.Item(9,8).Value = "NXT|Rup"
$ws.Cells.Item(9,9).Value = "NXT|Arthur"
$ws.Cells.Item(9,10).Value = "NXT|amos"
$ws.Cells.Item(9,11).Value = "FUT|GeRo"
$ws.Cells.Item(9,12).Value = "FUT|Nowy297"
$ws.Cells.Item(9,13).Value = "FUT|MeOw"
$ws.Cells.Item(9,14).Value = "20250723T161218.000Z"

# ---- Crystal Arcade (sheet index 2) ----
$ws = $wb.Worksheets.Item(2)

# Copy "Equipo 2" row formatting (style pattern from existing row 6)
$ws.Range("A6:N6").Copy()
$ws.Range("A18:N18").PasteSpecial(-4122)
$ws.Range("A20:N20").PasteSpecial(-4122)
$ws.Range("A25:N25").PasteSpecial(-4122)
$ws.Range("A26:N26").PasteSpecial(-4122)

# Copy "Equipo 1" row formatting (style pattern from existing row 4)
$ws.Range("A4:N4").Copy()
$ws.Range("A19:N19").PasteSpecial(-4122)
$ws.Range("A21:N21").PasteSpecial(-4122)
$ws.Range("A22:N22").PasteSpecial(-4122)
$ws.Range("A23:N23").PasteSpecial(-4122)
$ws.Range("A24:N24").PasteSpecial(-4122)

# Set cell values for new rows
$ws.Cells.Item(18,1).Value = "MAX"
$ws.Cells.Item(18,2).Value = "CORDELIUS"
$ws.Cells.Item(18,3).Value = "HANK"
$ws.Cells.Item(18,4).Value = "GALE"
$ws.Cells.Item(18,5).Value = "MEG"
$ws.Cells.Item(18,6).Value = "MR. P"
$ws.Cells.Item(18,7).Value = "Equipo 2"
$ws.Cells.Item(18,8).Value = "IC|Mebius"
$ws.Cells.Item(18,9).Value = "IC|Nob"
$ws.Cells.Item(18,10).Value = "IC|RamaZR"
$ws.Cells.Item(18,11).Value = "TTM|Angelboy"
$ws.Cells.Item(18,12).Value = "TTM|Maru"
$ws.Cells.Item(18,13).Value = "TTM|Maury"
$ws.Cells.Item(18,14).Value = "20250723T161006.000Z"

$ws.Cells.Item(19,1).Value = "MAX"
$ws.Cells.Item(19,2).Value = "CORDELIUS"
$ws.Cells.Item(19,3).Value = "HANK"
$ws.Cells.Item(19,4).Value = "GALE"
$ws.Cells.Item(19,5).Value = "MEG"
$ws.Cells.Item(19,6).Value = "MR. P"
$ws.Cells.Item(19,7).Value = "Equipo 1"
$ws.Cells.Item(19,8).Value = "IC|Mebius"
$ws.Cells.Item(19,9).Value = "IC|Nob"
$ws.Cells.Item(19,10).Value = "IC|RamaZR"
$ws.Cells.Item(19,11).Value = "TTM|Angelboy"
$ws.Cells.Item(19,12).Value = "TTM|Maru"
$ws.Cells.Item(19,13).Value = "TTM|Maury"
$ws.Cells.Item(19,14).Value = "20250723T160730.000Z"

$ws.Cells.Item(20,1).Value = "MAX"
$ws.Cells.Item(20,2).Value = "CORDELIUS"
$ws.Cells.Item(20,3).Value = "HANK"
$ws.Cells.Item(20,4).Value = "GALE"
$ws.Cells.Item(20,5).Value = "MEG"
$ws.Cells.Item(20,6).Value = "MR. P"
$ws.Cells.Item(20,7).Value = "Equipo 2"
$ws.Cells.Item(20,8).Value = "IC|Mebius"
$ws.Cells.Item(20,9).Value = "IC|Nob"
$ws.Cells.Item(20,10).Value = "IC|RamaZR"
$ws.Cells.Item(20,11).Value = "TTM|Angelboy"
$ws.Cells.Item(20,12).Value = "TTM|Maru"
$ws.Cells.Item(20,13).Value = "TTM|Maury"
$ws.Cells.Item(20,14).Value = "20250723T160543.000Z"

$ws.Cells.Item(21,1).Value = "AMBER"
$ws.Cells.Item(21,2).Value = "MEEPLE"
$ws.Cells.Item(21,3).Value = "HANK"
$ws.Cells.Item(21,4).Value = "KENJI"
$ws.Cells.Item(21,5).Value = "CORDELIUS"
$ws.Cells.Item(21,6).Value = "MR. P"
$ws.Cells.Item(21,7).Value = "Equipo 1"
$ws.Cells.Item(21,8).Value = "IC|Mebius"
$ws.Cells.Item(21,9).Value = "IC|Nob"
$ws.Cells.Item(21,10).Value = "IC|RamaZR"
$ws.Cells.Item(21,11).Value = "TTM|Angelboy"
$ws.Cells.Item(21,12).Value = "TTM|Maru"
$ws.Cells.Item(21,13).Value = "TTM|Maury"
$ws.Cells.Item(21,14).Value = "20250723T160005.000Z"

$ws.Cells.Item(22,1).Value = "AMBER"
$ws.Cells.Item(22,2).Value = "MEEPLE"
$ws.Cells.Item(22,3).Value = "HANK"
$ws.Cells.Item(22,4).Value = "KENJI"
$ws.Cells.Item(22,5).Value = "CORDELIUS"
$ws.Cells.Item(22,6).Value = "MR. P"
$ws.Cells.Item(22,7).Value = "Equipo 1"
$ws.Cells.Item(22,8).Value = "IC|Mebius"
$ws.Cells.Item(22,9).Value = "IC|Nob"
$ws.Cells.Item(22,10).Value = "IC|RamaZR"
$ws.Cells.Item(22,11).Value = "TTM|Angelboy"
$ws.Cells.Item(22,12).Value = "TTM|Maru"
$ws.Cells.Item(22,13).Value = "TTM|Maury"
$ws.Cells.Item(22,14).Value = "20250723T155807.000Z"

$ws.Cells.Item(23,1).Value = "KENJI"
$ws.Cells.Item(23,2).Value = "MR. P"
$ws.Cells.Item(23,3).Value = "CORDELIUS"
$ws.Cells.Item(23,4).Value = "ALLI"
$ws.Cells.Item(23,5).Value = "SPIKE"
$ws.Cells.Item(23,6).Value = "FINX"
$ws.Cells.Item(23,7).Value = "Equipo 1"
$ws.Cells.Item(23,8).Value = "HMB|BosS"
$ws.Cells.Item(23,9).Value = "HMB|Lukii"
$ws.Cells.Item(23,10).Value = "HMB|Symantec"
$ws.Cells.Item(23,11).Value = "TH|LeNain"
$ws.Cells.Item(23,12).Value = "TH|iKaoss"
$ws.Cells.Item(23,13).Value = "TH|Zhar"
$ws.Cells.Item(23,14).Value = "20250723T162040.000Z"

$ws.Cells.Item(24,1).Value = "KENJI"
$ws.Cells.Item(24,2).Value = "MR. P"
$ws.Cells.Item(24,3).Value = "CORDELIUS"
$ws.Cells.Item(24,4).Value = "ALLI"
$ws.Cells.Item(24,5).Value = "SPIKE"
$ws.Cells.Item(24,6).Value = "FINX"
$ws.Cells.Item(24,7).Value = "Equipo 1"
$ws.Cells.Item(24,8).Value = "HMB|BosS"
$ws.Cells.Item(24,9).Value = "HMB|Lukii"
$ws.Cells.Item(24,10).Value = "HMB|Symantec"
$ws.Cells.Item(24,11).Value = "TH|LeNain"
$ws.Cells.Item(24,12).Value = "TH|iKaoss"
$ws.Cells.Item(24,13).Value = "TH|Zhar"
$ws.Cells.Item(24,14).Value = "20250723T161806.000Z"

$ws.Cells.Item(25,1).Value = "ALLI"
$ws.Cells.Item(25,2).Value = "TARA"
$ws.Cells.Item(25,3).Value = "JANET"
$ws.Cells.Item(25,4).Value = "LUMI"
$ws.Cells.Item(25,5).Value = "GUS"
$ws.Cells.Item(25,6).Value = "LILY"
$ws.Cells.Item(25,7).Value = "Equipo 2"
$ws.Cells.Item(25,8).Value = "HMB|Lukii"
$ws.Cells.Item(25,9).Value = "HMB|BosS"
$ws.Cells.Item(25,10).Value = "HMB|Symantec"
$ws.Cells.Item(25,11).Value = "TH|iKaoss"
$ws.Cells.Item(25,12).Value = "TH|Zhar"
$ws.Cells.Item(25,13).Value = "TH|LeNain"
$ws.Cells.Item(25,14).Value = "20250723T161133.000Z"

$ws.Cells.Item(26,1).Value = "ALLI"
$ws.Cells.Item(26,2).Value = "TARA"
$ws.Cells.Item(26,3).Value = "JANET"
$ws.Cells.Item(26,4).Value = "LUMI"
$ws.Cells.Item(26,5).Value = "GUS"
$ws.Cells.Item(26,6).Value = "LILY"
$ws.Cells.Item(26,7).Value = "Equipo 2"
$ws.Cells.Item(26,8).Value = "HMB|Lukii"
$ws.Cells.Item(26,9).Value = "HMB|BosS"
$ws.Cells.Item(26,10).Value = "HMB|Symantec"
$ws.Cells.Item(26,11).Value = "TH|iKaoss"
$ws.Cells.Item(26,12).Value = "TH|Zhar"
$ws.Cells.Item(26,13).Value = "TH|LeNain"
$ws.Cells.Item(26,14).Value = "20250723T160930.000Z"

# ---- Dueling Beetles (sheet index 9) ----
$ws = $wb.Worksheets.Item(9)

# Copy "Equipo 1" row formatting (style pattern from existing row 4)
$ws.Range("A4:N4").Copy()
$ws.Range("A17:N17").PasteSpecial(-4122)

# Set cell values for new rows
$ws.Cells.Item(17,1).Value = "CORDELIUS"
$ws.Cells.Item(17,2).Value = "HANK"
$ws.Cells.Item(17,3).Value = "GRAY"
$ws.Cells.Item(17,4).Value = "DRACO"
$ws.Cells.Item(17,5).Value = "LUMI"
$ws.Cells.Item(17,6).Value = "CHARLIE"
$ws.Cells.Item(17,7).Value = "Equipo 1"
$ws.Cells.Item(17,8).Value = "HMB|BosS"
$ws.Cells.Item(17,9).Value = "HMB|Symantec"
$ws.Cells.Item(17,10).Value = "HMB|Lukii"
$ws.Cells.Item(17,11).Value = "TH|LeNain"
$ws.Cells.Item(17,12).Value = "TH|iKaoss"
$ws.Cells.Item(17,13).Value = "TH|Zhar"
$ws.Cells.Item(17,14).Value = "20250723T162626.000Z"

